$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 310
$ws.Range("F4").Value = 3162
$ws.Range("F9").Value = 7409
$ws.Range("F12").Value = 451
$ws.Range("F14").Value = 1818
$ws.Range("F15").Value = 1101
$ws.Range("F16").Value = 41
$ws.Range("F18").Value = 1872
$ws.Range("F19").Value = 1397
$ws.Range("F20").Value = 1278
$ws.Range("F21").Value = 672
$ws.Range("F23").Value = 1165
$ws.Range("F24").Value = 128
$ws.Range("F25").Value = 561
$ws.Range("F27").Value = 139
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 4760
$ws.Range("F30").Value = 2292
$ws.Range("F31").Value = 3960
$ws.Range("F32").Value = 2211
$ws.Range("F33").Value = 181
$ws.Range("F34").Value = 231
$ws.Range("F35").Value = 1158
$ws.Range("F38").Value = 56
$ws.Range("F39").Value = 391
$ws.Range("F42").Value = 537
$ws.Range("F43").Value = 277
$ws.Range("F44").Value = 205
$ws.Range("F45").Value = 841
$ws.Range("F46").Value = 442
$ws.Range("F47").Value = 38
$ws.Range("F49").Value = 178

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 166
$ws.Range("F16").Value = 558
$ws.Range("F25").Value = 106
$ws.Range("F29").Value = 12
$ws.Range("F32").Value = 1358
$ws.Range("F33").Value = 1358
$ws.Range("F35").Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1870
$ws.Range("F9").Value = 1156
$ws.Range("F10").Value = 1158
$ws.Range("F12").Value = 462
$ws.Range("F13").Value = 1876
$ws.Range("F14").Value = 8285
$ws.Range("F15").Value = 550

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 310
$ws.Range("F4").Value = 3162
$ws.Range("F5").Value = 1870
$ws.Range("F8").Value = 7409
$ws.Range("F9").Value = 1156
$ws.Range("F10").Value = 1158
$ws.Range("F12").Value = 462
$ws.Range("F13").Value = 451
$ws.Range("F16").Value = 1101
$ws.Range("F17").Value = 41
$ws.Range("F19").Value = 1872
$ws.Range("F20").Value = 1397
$ws.Range("F21").Value = 1278
$ws.Range("F23").Value = 672
$ws.Range("F25").Value = 1165
$ws.Range("F27").Value = 128
$ws.Range("F29").Value = 558
$ws.Range("F30").Value = 561
$ws.Range("F33").Value = 139
$ws.Range("F34").Value = 83
$ws.Range("F35").Value = 4760
$ws.Range("F36").Value = 2292
$ws.Range("F37").Value = 3960
$ws.Range("F38").Value = 181
$ws.Range("F39").Value = 231
$ws.Range("F40").Value = 1158
$ws.Range("F42").Value = 56
$ws.Range("F44").Value = 106
$ws.Range("F45").Value = 537
$ws.Range("F46").Value = 277
$ws.Range("F48").Value = 442
$ws.Range("F49").Value = 1359
